$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1624.303
$ws.Range("I137").Value = 1026.9714
$ws.Range("J137").Value = 2298.7097
$ws.Range("K137").Value = 3080.9142
$ws.Range("L137").Value = 6896.1291
$ws.Range("M137").Value = -530.9141999999997
$ws.Range("N137").Value = -11996.1291

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17135.238
$ws.Range("I32").Value = 16323.204
$ws.Range("K32").Value = 16323.204
$ws.Range("M32").Value = -16036.204

$ws.Range("H61").Value = 2393.2727
$ws.Range("I61").Value = 2192.2104
$ws.Range("J61").Value = 3666.6667
$ws.Range("K61").Value = 2192.2104
$ws.Range("L61").Value = 3666.6667
$ws.Range("M61").Value = -1980.2104
$ws.Range("N61").Value = -4090.6667

$ws.Range("H74").Value = 1424.4681
$ws.Range("I74").Value = 1292.6842
$ws.Range("J74").Value = 1980.8889
$ws.Range("K74").Value = 1292.6842
$ws.Range("L74").Value = 1980.8889
$ws.Range("M74").Value = -418.6841999999999
$ws.Range("N74").Value = -3728.8889

$ws.Range("H77").Value = 1424.4681
$ws.Range("I77").Value = 1292.6842
$ws.Range("J77").Value = 1980.8889
$ws.Range("K77").Value = 6463.420999999999
$ws.Range("L77").Value = 9904.4445
$ws.Range("M77").Value = -2095.420999999999
$ws.Range("N77").Value = -18640.4445

$ws.Range("H104").Value = 49225
$ws.Range("J104").Value = 49225
$ws.Range("L104").Value = 49225
$ws.Range("N104").Value = -56213

$ws.Range("H136").Value = 2393.2727
$ws.Range("I136").Value = 2192.2104
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 6576.6312
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -4026.6312
$ws.Range("N136").Value = -16100.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 70357.07000000001
$ws.Range("I86").Value = 3680.889
$ws.Range("J86").Value = 170371.33
$ws.Range("K86").Value = 3680.889
$ws.Range("L86").Value = 170371.33
$ws.Range("M86").Value = -2557.889
$ws.Range("N86").Value = -172617.33

$ws.Range("H89").Value = 70357.07000000001
$ws.Range("I89").Value = 3680.889
$ws.Range("J89").Value = 170371.33
$ws.Range("K89").Value = 18404.445
$ws.Range("L89").Value = 851856.6499999999
$ws.Range("M89").Value = -12788.445
$ws.Range("N89").Value = -863088.6499999999

$ws.Range("H134").Value = 2270.4333
$ws.Range("I134").Value = 1777.8636
$ws.Range("J134").Value = 3625
$ws.Range("K134").Value = 5333.5908
$ws.Range("L134").Value = 10875
$ws.Range("M134").Value = -2798.5908
$ws.Range("N134").Value = -15945

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1613.8276
$ws.Range("I31").Value = 1830.9269
$ws.Range("J31").Value = 1420.326
$ws.Range("K31").Value = 1830.9269
$ws.Range("L31").Value = 1420.326
$ws.Range("M31").Value = -1535.9269
$ws.Range("N31").Value = -2010.326

$ws.Range("H34").Value = 1613.8276
$ws.Range("I34").Value = 1830.9269
$ws.Range("J34").Value = 1420.326
$ws.Range("K34").Value = 1830.9269
$ws.Range("L34").Value = 1420.326
$ws.Range("M34").Value = -1628.9269
$ws.Range("N34").Value = -1824.326

$ws.Range("H58").Value = 1324378.1
$ws.Range("I58").Value = 1950940.8
$ws.Range("K58").Value = 1950940.8
$ws.Range("M58").Value = -1950737.8

$ws.Range("H132").Value = 798301.5600000001
$ws.Range("I132").Value = 2704487.8
$ws.Range("J132").Value = 4057.3333
$ws.Range("K132").Value = 8113463.399999999
$ws.Range("L132").Value = 12171.9999
$ws.Range("M132").Value = -8110933.399999999
$ws.Range("N132").Value = -17231.9999

$ws.Range("H134").Value = 1424.9231
$ws.Range("I134").Value = 1210
$ws.Range("J134").Value = 2607
$ws.Range("K134").Value = 3630
$ws.Range("L134").Value = 7821
$ws.Range("M134").Value = -1095
$ws.Range("N134").Value = -12891

$ws.Range("H136").Value = 1324378.1
$ws.Range("I136").Value = 1950940.8
$ws.Range("K136").Value = 5852822.4
$ws.Range("M136").Value = -5850272.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1023.8461
$ws.Range("J34").Value = 2950
$ws.Range("L34").Value = 8850
$ws.Range("N34").Value = -9018

$ws.Range("H39").Value = 4213.0835
$ws.Range("J39").Value = 4925.7
$ws.Range("L39").Value = 14777.1
$ws.Range("N39").Value = -15365.1

$ws.Range("H48").Value = 6000
$ws.Range("J48").Value = 6000
$ws.Range("L48").Value = 18000
$ws.Range("N48").Value = -18500

$ws.Range("H54").Value = 6666.6665
$ws.Range("J54").Value = 6666.6665
$ws.Range("L54").Value = 19999.9995
$ws.Range("N54").Value = -21117.9995

$ws.Range("H59").Value = 5000
$ws.Range("J59").Value = 5000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16080

$ws.Range("H63").Value = 236491.56
$ws.Range("I63").Value = 352404
$ws.Range("K63").Value = 1057212
$ws.Range("M63").Value = -1056463

$ws.Range("H66").Value = 236491.56
$ws.Range("I66").Value = 352404
$ws.Range("K66").Value = 3171636
$ws.Range("M66").Value = -3167892

$ws.Range("H68").Value = 1282.1266
$ws.Range("J68").Value = 1642.2368
$ws.Range("L68").Value = 4926.7104
$ws.Range("N68").Value = -6548.7104

$ws.Range("H71").Value = 1282.1266
$ws.Range("J71").Value = 1642.2368
$ws.Range("L71").Value = 14780.1312
$ws.Range("N71").Value = -22892.1312

$ws.Range("H75").Value = 8551.75
$ws.Range("J75").Value = 10912.777
$ws.Range("L75").Value = 32738.331
$ws.Range("N75").Value = -34734.331

$ws.Range("H78").Value = 8551.75
$ws.Range("J78").Value = 10912.777
$ws.Range("L78").Value = 98214.993
$ws.Range("N78").Value = -108198.993

$ws.Range("H94").Value = 4262.857
$ws.Range("J94").Value = 4656.6665
$ws.Range("L94").Value = 13969.9995
$ws.Range("N94").Value = -15321.9995

$ws.Range("H96").Value = 32012.4
$ws.Range("I96").Value = 25512.5
$ws.Range("J96").Value = 36345.668
$ws.Range("K96").Value = 76537.5
$ws.Range("L96").Value = 109037.004
$ws.Range("M96").Value = -74478.5
$ws.Range("N96").Value = -113155.004

$ws.Range("H102").Value = 7309.6895
$ws.Range("J102").Value = 7309.6895
$ws.Range("L102").Value = 21929.0685
$ws.Range("N102").Value = -26797.0685

$ws.Range("H109").Value = 1500
$ws.Range("I109").Value = 1500
$ws.Range("K109").Value = 4500
$ws.Range("M109").Value = -3460

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = 0

$ws.Range("H119").Value = 7666.6665
$ws.Range("I119").Value = 4000
$ws.Range("J119").Value = 15000
$ws.Range("K119").Value = 12000
$ws.Range("L119").Value = 45000
$ws.Range("M119").Value = -7162
$ws.Range("N119").Value = -54676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 50000
$ws.Range("I33").Value = 50000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 50000
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -49748

$ws.Range("H126").Value = 1985.7391
$ws.Range("I126").Value = 1269.8
$ws.Range("J126").Value = 2536.4614
$ws.Range("K126").Value = 3809.4
$ws.Range("L126").Value = 7609.3842
$ws.Range("M126").Value = -1339.4
$ws.Range("N126").Value = -12549.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I68").Value = 1768.0667
$ws.Range("J68").Value = 3667.1667
$ws.Range("K68").Value = 1768.0667
$ws.Range("L68").Value = 3667.1667
$ws.Range("M68").Value = -1019.0667
$ws.Range("N68").Value = -5165.1667

$ws.Range("I71").Value = 1768.0667
$ws.Range("J71").Value = 3667.1667
$ws.Range("K71").Value = 8840.333500000001
$ws.Range("L71").Value = 18335.8335
$ws.Range("M71").Value = -5096.333500000001
$ws.Range("N71").Value = -25823.8335

$ws.Range("H132").Value = 5456.92
$ws.Range("I132").Value = 5456.92
$ws.Range("K132").Value = 16370.76
$ws.Range("M132").Value = -13840.76

$ws.Range("H136").Value = 1990.9445
$ws.Range("I136").Value = 1346.96
$ws.Range("K136").Value = 4040.88
$ws.Range("M136").Value = -1490.88

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2188.4358
$ws.Range("I132").Value = 2094.423
$ws.Range("J132").Value = 2376.4614
$ws.Range("K132").Value = 6283.268999999999
$ws.Range("L132").Value = 7129.3842
$ws.Range("M132").Value = -3753.268999999999
$ws.Range("N132").Value = -12189.3842

$ws.Range("H136").Value = 2622.7144
$ws.Range("I136").Value = 2336.182
$ws.Range("J136").Value = 3107.6155
$ws.Range("K136").Value = 7008.545999999999
$ws.Range("L136").Value = 9322.8465
$ws.Range("M136").Value = -4458.545999999999
$ws.Range("N136").Value = -14422.8465
